# Add "Norway" and "Poland" market test-data sheets.
#
# New sheets are cloned from "Turkey" (a sheet whose row heights are all
# default and whose column D is already the 21.109375-wide, non-bestFit
# layout used by the new country sheets - unlike "Hungary", whose rows
# 3-5 carry an explicit ht="28.8" and whose column D is bestFit at
# 8.44140625) so the new sheets come out with the same layout as in the
# target workbook.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Turkey")
$hungary  = $wb.Worksheets.Item("Hungary")

# --- Norway: inserted right after Hungary -----------------------------
$template.Copy($null, $hungary)
$norway = $wb.Worksheets.Item($hungary.Index + 1)
$norway.Name = "Norway"

# Set B4 (ticket code) before B2 (market name) so the new shared strings
# land in the same order as the target file (NGC code, then "X Market").
$norway.Range("B4").Value = "NGC-2931/T3068"
$norway.Range("B2").Value = "Norway Market"
$norway.Range("H22").Select()

# --- Poland: inserted right after Norway -------------------------------
$template.Copy($null, $norway)
$poland = $wb.Worksheets.Item($norway.Index + 1)
$poland.Name = "Poland"

$poland.Range("B4").Value = "NGC-2920/3034"
$poland.Range("B2").Value = "Poland Market"
$poland.Range("H22").Select()

# Norway is the active/selected tab in the saved workbook (Poland, being
# the most recently inserted sheet, would otherwise keep that status).
$norway.Activate()
